# Scheduled-runner style update: refresh market-board derived price/profit
# columns (H:N) for a handful of Leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1430.2122
$ws.Range("I28").Value = 229.47826
$ws.Range("J28").Value = 4191.9
$ws.Range("K28").Value = 229.47826
$ws.Range("L28").Value = 4191.9
$ws.Range("M28").Value = 255.52174
$ws.Range("N28").Value = -5161.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4816.6665
$ws.Range("I116").Value = 1950
$ws.Range("J116").Value = 6250
$ws.Range("K116").Value = 1950
$ws.Range("L116").Value = 6250
$ws.Range("M116").Value = 1492
$ws.Range("N116").Value = -13134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35715124
$ws.Range("J2").Value = 1193.1
$ws.Range("L2").Value = 1193.1
$ws.Range("N2").Value = -1419.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1464.08
$ws.Range("I32").Value = 1464.08
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1464.08
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1177.08
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 166667420
$ws.Range("I45").Value = 166667420
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 166667420
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -166667043
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1807.9131
$ws.Range("I61").Value = 1346.421
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1346.421
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1134.421
$ws.Range("N61").Value = -4424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 812.1475
$ws.Range("I74").Value = 720.1321
$ws.Range("J74").Value = 1421.75
$ws.Range("K74").Value = 720.1321
$ws.Range("L74").Value = 1421.75
$ws.Range("M74").Value = 153.8679
$ws.Range("N74").Value = -3169.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 812.1475
$ws.Range("I77").Value = 720.1321
$ws.Range("J77").Value = 1421.75
$ws.Range("K77").Value = 3600.6605
$ws.Range("L77").Value = 7108.75
$ws.Range("M77").Value = 767.3395
$ws.Range("N77").Value = -15844.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 35715124
$ws.Range("J116").Value = 1193.1
$ws.Range("L116").Value = 1193.1
$ws.Range("N116").Value = -5781.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1579.8
$ws.Range("I122").Value = 1571.1818
$ws.Range("J122").Value = 1603.5
$ws.Range("K122").Value = 4713.5454
$ws.Range("L122").Value = 4810.5
$ws.Range("M122").Value = -2263.5454
$ws.Range("N122").Value = -9710.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3940
$ws.Range("I132").Value = 3873.8958
$ws.Range("J132").Value = 4468.8335
$ws.Range("K132").Value = 11621.6874
$ws.Range("L132").Value = 13406.5005
$ws.Range("M132").Value = -9091.687399999999
$ws.Range("N132").Value = -18466.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1807.9131
$ws.Range("I136").Value = 1346.421
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 4039.263
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1489.263
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35715124
$ws.Range("J3").Value = 1193.1
$ws.Range("L3").Value = 1193.1
$ws.Range("N3").Value = -1421.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2309.2246
$ws.Range("I105").Value = 1835.8966
$ws.Range("J105").Value = 2995.55
$ws.Range("K105").Value = 1835.8966
$ws.Range("L105").Value = 2995.55
$ws.Range("M105").Value = -88.89660000000003
$ws.Range("N105").Value = -6489.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1276.4193
$ws.Range("I107").Value = 1264.3478
$ws.Range("J107").Value = 1311.125
$ws.Range("K107").Value = 1264.3478
$ws.Range("L107").Value = 1311.125
$ws.Range("M107").Value = 655.6522
$ws.Range("N107").Value = -5151.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2691021.2
$ws.Range("I31").Value = 2098.1191
$ws.Range("J31").Value = 8337760
$ws.Range("K31").Value = 2098.1191
$ws.Range("L31").Value = 8337760
$ws.Range("M31").Value = -1803.1191
$ws.Range("N31").Value = -8338350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2691021.2
$ws.Range("I34").Value = 2098.1191
$ws.Range("J34").Value = 8337760
$ws.Range("K34").Value = 2098.1191
$ws.Range("L34").Value = 8337760
$ws.Range("M34").Value = -1896.1191
$ws.Range("N34").Value = -8338164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 14706986
$ws.Range("I113").Value = 31250744
$ws.Range("J113").Value = 1422.2222
$ws.Range("K113").Value = 31250744
$ws.Range("L113").Value = 1422.2222
$ws.Range("M113").Value = -31248574
$ws.Range("N113").Value = -5762.2222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1519.5946
$ws.Range("J126").Value = 1491.4
$ws.Range("L126").Value = 4474.200000000001
$ws.Range("N126").Value = -9414.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 79854.305
$ws.Range("I132").Value = 126891
$ws.Range("J132").Value = 4595.6
$ws.Range("K132").Value = 380673
$ws.Range("L132").Value = 13786.8
$ws.Range("M132").Value = -378143
$ws.Range("N132").Value = -18846.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 49365.855
$ws.Range("J102").Value = 49365.855
$ws.Range("L102").Value = 49365.855
$ws.Range("N102").Value = -55855.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5093.143
$ws.Range("I122").Value = 7614.857
$ws.Range("J122").Value = 2571.4285
$ws.Range("K122").Value = 22844.571
$ws.Range("L122").Value = 7714.2855
$ws.Range("M122").Value = -20394.571
$ws.Range("N122").Value = -12614.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5724.347
$ws.Range("I132").Value = 7514.3125
$ws.Range("J132").Value = 2355
$ws.Range("K132").Value = 22542.9375
$ws.Range("L132").Value = 7065
$ws.Range("M132").Value = -20012.9375
$ws.Range("N132").Value = -12125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2703.2222
$ws.Range("I122").Value = 2534.15
$ws.Range("J122").Value = 3186.2856
$ws.Range("K122").Value = 7602.450000000001
$ws.Range("L122").Value = 9558.856800000001
$ws.Range("M122").Value = -5152.450000000001
$ws.Range("N122").Value = -14458.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1040.5312
$ws.Range("I136").Value = 963.6875
$ws.Range("K136").Value = 2891.0625
$ws.Range("M136").Value = -341.0625
